# Insert a new header row of numeric column indices (0-11, 12) above the
# existing data, shifting the current header row (with text labels like
# "Lg.", "Threading", ...) and all data rows down by one.
#
# Before: row 1 = text headers, rows 2-58 = data
# After:  row 1 = numeric indices 0..12, row 2 = text headers, rows 3-59 = data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at the very top; everything else shifts down by one.
$ws.Rows.Item(1).Insert()

# The newly inserted row 1 has no formatting. Copy the formatting (bold font,
# borders, centered alignment) that used to belong to the header row (now
# row 2) onto the new row 1.
$ws.Range("A2:M2").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# The old header row (now row 2) should no longer carry that special
# formatting.
$ws.Range("A2:M2").ClearFormats()

# Fill the new row 1 with the numeric sequence 0..12.
for ($i = 0; $i -lt 13; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $i
}

# In the shifted header row (row 2), columns J, L and M end up blank.
$ws.Range("J2").ClearContents()
$ws.Range("L2").ClearContents()
$ws.Range("M2").ClearContents()
